# PROVA - spostate commesse esterne non tassative al gruppo 2
# Rewrite the data rows on the "Release Date (RD)" sheet (rows 2-13),
# adding the new "commesse esterne non tassative" entries while keeping
# the header row (row 1) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Date (RD)")

$data = @(
    @(253371, 45910.58333333334, 0),
    @(253367, 45911.58333333334, 0),
    @(253472, 45912.58333333334, 0),
    @(253549, 45912.58333333334, 0),
    @(253392, 45911.58333333334, 0),
    @(253393, 45912.58333333334, 0),
    @(253295, 45911.58333333334, "X"),
    @(252397, 45911.58333333334, 0),
    @(253244, 45912.58333333334, 0),
    @(252274, 45911.58333333334, 0),
    @(253668, 45910.58333333334, 0),
    @(252741, 45911.58333333334, 0)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 2).NumberFormat = "yyyy-mm-dd h:mm:ss"
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
